# Insert 3 new rows at row 865 (shifts all existing rows 865:930 down to 868:933)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A865:A867").EntireRow.Insert()

# Row 865 - new weekly entry, Calidad "Primera"
$ws.Cells.Item(865,1).Value  = 6
$ws.Cells.Item(865,2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(865,3).Value  = "Metropolitana"
$ws.Cells.Item(865,4).Value  = 44714
$ws.Cells.Item(865,5).Value  = 13
$ws.Cells.Item(865,6).Value  = 100114014
$ws.Cells.Item(865,7).Value  = "Betarraga"
$ws.Cells.Item(865,8).Value  = "Sin especificar"
$ws.Cells.Item(865,9).Value  = "Primera"
$ws.Cells.Item(865,10).Value = 46000
$ws.Cells.Item(865,11).Value = 100
$ws.Cells.Item(865,12).Value = 110
$ws.Cells.Item(865,13).Value = 105
$ws.Cells.Item(865,14).Value = "`$/unidad"
$ws.Cells.Item(865,15).Value = "Región Metropolitana"
$ws.Cells.Item(865,16).Value = 105
$ws.Cells.Item(865,17).Value = 1
$ws.Cells.Item(865,18).Value = "Hortaliza"

# Row 866 - new weekly entry, Calidad "Segunda"
$ws.Cells.Item(866,1).Value  = 6
$ws.Cells.Item(866,2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(866,3).Value  = "Metropolitana"
$ws.Cells.Item(866,4).Value  = 44714
$ws.Cells.Item(866,5).Value  = 13
$ws.Cells.Item(866,6).Value  = 100114014
$ws.Cells.Item(866,7).Value  = "Betarraga"
$ws.Cells.Item(866,8).Value  = "Sin especificar"
$ws.Cells.Item(866,9).Value  = "Segunda"
$ws.Cells.Item(866,10).Value = 39000
$ws.Cells.Item(866,11).Value = 80
$ws.Cells.Item(866,12).Value = 85
$ws.Cells.Item(866,13).Value = 82
$ws.Cells.Item(866,14).Value = "`$/unidad"
$ws.Cells.Item(866,15).Value = "Región Metropolitana"
$ws.Cells.Item(866,16).Value = 82
$ws.Cells.Item(866,17).Value = 1
$ws.Cells.Item(866,18).Value = "Hortaliza"

# Row 867 - new weekly entry, Calidad "Tercera"
$ws.Cells.Item(867,1).Value  = 6
$ws.Cells.Item(867,2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(867,3).Value  = "Metropolitana"
$ws.Cells.Item(867,4).Value  = 44714
$ws.Cells.Item(867,5).Value  = 13
$ws.Cells.Item(867,6).Value  = 100114014
$ws.Cells.Item(867,7).Value  = "Betarraga"
$ws.Cells.Item(867,8).Value  = "Sin especificar"
$ws.Cells.Item(867,9).Value  = "Tercera"
$ws.Cells.Item(867,10).Value = 9000
$ws.Cells.Item(867,11).Value = 65
$ws.Cells.Item(867,12).Value = 65
$ws.Cells.Item(867,13).Value = 65
$ws.Cells.Item(867,14).Value = "`$/unidad"
$ws.Cells.Item(867,15).Value = "Región Metropolitana"
$ws.Cells.Item(867,16).Value = 65
$ws.Cells.Item(867,17).Value = 1
$ws.Cells.Item(867,18).Value = "Hortaliza"

# Make sure the date cells use the same date/time number format as the rest of column D
$ws.Cells.Item(865,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(866,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(867,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
